# Updated literature IPM list
# Adds the new Madin, Hughes & Connolly (2012) PLoS One reference on the
# "Original reference list" sheet, and reflects the corresponding updates
# on the "Species statistics" sheet (pivot summary row + chart source cell).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. "Original reference list" sheet: insert the new publication as row 47,
#    pushing the existing row 47 ("Metcalf & Mitchell-Olds" ...) and
#    everything below it down by one row.
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Original reference list")

$ws1.Rows.Item(47).Insert()

$ws1.Cells.Item(47, 1).Value = "Madin, Hughes & Connolly"
$ws1.Cells.Item(47, 2).Value = "PLoS One"
$ws1.Cells.Item(47, 3).Value = 2012
$ws1.Cells.Item(47, 4).Value = "Animalia"
$ws1.Cells.Item(47, 5).Value = "Acropora hyacinthus"
$ws1.Cells.Item(47, 7).Value = "Madin JS, Hughes TP & Connolly SR (2012) Calcification, storm damage and population resilience of tabular corals under climate change. PLoS One 7: 1-10"

# ---------------------------------------------------------------------------
# 2. "Species statistics" sheet: the new reference brings the 2012
#    publication count (chart source table) from 24 to 25, and adds a new
#    "Acropora hyacinthus" line (count 1) to the Count-of-Species pivot
#    summary just above the Grand Total row, bumping the grand total from
#    98 to 99.
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Species statistics")

$ws2.Cells.Item(36, 6).Value = 25

$ws2.Rows.Item(65).Insert()
$ws2.Cells.Item(65, 1).Value = "Acropora hyacinthus"
$ws2.Cells.Item(65, 2).Value = 1
$ws2.Cells.Item(66, 2).Value = 99

# ---------------------------------------------------------------------------
# 3. Cosmetic window / selection state (best effort).
# ---------------------------------------------------------------------------
try {
    $excel.ActiveWindow.Left = 27320
    $excel.ActiveWindow.Top = -5900
    $excel.ActiveWindow.Width = 28000
    $excel.ActiveWindow.Height = 21000
} catch {
}

$ws2.Activate()
$ws2.Range("G26").Select()
